$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.576.75"
$ws.Range("E2").Value = "  +2.14%  "

# Row 3
$ws.Range("D3").Value = "2.083.70"
$ws.Range("E3").Value = "  +10.19%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.48"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6
$ws.Range("E6").Value = "  -4.15%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.09"
$ws.Range("E8").Value = "  +4.58%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.85"
$ws.Range("E9").Value = "  +7.80%  "

# Row 10
$ws.Range("E10").Value = "  +1.98%  "

# Row 11
$ws.Range("E11").Value = "  -4.17%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0988"
$ws.Range("E12").Value = "  +0.27%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.50"
$ws.Range("E13").Value = "  -3.12%  "

# Row 14
$ws.Range("D14").Value = "2.386.64"
$ws.Range("E14").Value = "  +10.17%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.826"
$ws.Range("E15").Value = "  +4.37%  "

# Row 16
$ws.Range("D16").Value = "2.074.92"
$ws.Range("E16").Value = "  +9.62%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.91"
$ws.Range("E17").Value = "  -2.92%  "

# Row 18
$ws.Range("D18").Value = "36.563.87"
$ws.Range("E18").Value = "  +2.24%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.67"
$ws.Range("E19").Value = "  -2.70%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0814"
$ws.Range("E20").Value = "  -2.21%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "238.56"
$ws.Range("E21").Value = "  -3.49%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.71"
$ws.Range("E22").Value = "  -2.99%  "

# Row 23
$ws.Range("E23").Value = "  -4.86%  "

# Row 24
$ws.Range("E24").Value = "  -0.03%  "

# Row 25
$ws.Range("E25").Value = "  -8.30%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.81"
$ws.Range("E26").Value = "  +1.63%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.39"
$ws.Range("E27").Value = "  +10.46%  "

# Row 28
$ws.Range("E28").Value = "  +2.27%  "

# Row 29
$ws.Range("E29").Value = "  -8.45%  "

# Row 30
$ws.Range("E30").Value = "  -5.28%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.27"
$ws.Range("E31").Value = "  +57.89%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.36"

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0583"
$ws.Range("E33").Value = "  -4.47%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0911"
$ws.Range("E34").Value = "  +17.57%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.90"
$ws.Range("E35").Value = "  +2.75%  "

# Row 36
$ws.Range("E36").Value = "  -0.07%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.29"
$ws.Range("E37").Value = "  +18.25%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.901"
$ws.Range("E38").Value = "  +5.18%  "

# Row 39
$ws.Range("E39").Value = "  -6.45%  "

# Row 40
$ws.Range("E40").Value = "  -9.25%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.14"
$ws.Range("E41").Value = "  +5.14%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.79"
$ws.Range("E42").Value = "  -1.42%  "

# Row 43
$ws.Range("E43").Value = "  -6.22%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.78"
$ws.Range("E44").Value = "  +15.52%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.91"
$ws.Range("E45").Value = "  -6.26%  "

# Row 46
$ws.Range("D46").Value = "1.331.39"
$ws.Range("E46").Value = "  +1.10%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.83"
$ws.Range("E48").Value = "  +4.00%  "

# Row 49
$ws.Range("D49").Value = "2.277.56"

# Row 50
$ws.Range("E50").Value = "  -5.64%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.86"
$ws.Range("E51").Value = "  +15.47%  "
